$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 883
$ws.Cells.Item(5, 6).Value = 1185
$ws.Cells.Item(6, 6).Value = 67
$ws.Cells.Item(7, 6).Value = 4346
$ws.Cells.Item(8, 6).Value = 2592
$ws.Cells.Item(10, 6).Value = 2505
$ws.Cells.Item(14, 6).Value = 1655
$ws.Cells.Item(15, 6).Value = 658
$ws.Cells.Item(16, 6).Value = 19
$ws.Cells.Item(18, 6).Value = 319
$ws.Cells.Item(20, 6).Value = 271
$ws.Cells.Item(21, 6).Value = 74
$ws.Cells.Item(22, 6).Value = 23
$ws.Cells.Item(23, 6).Value = 476
$ws.Cells.Item(26, 6).Value = 540
$ws.Cells.Item(28, 6).Value = 102
$ws.Cells.Item(30, 6).Value = 400
$ws.Cells.Item(31, 6).Value = 47
$ws.Cells.Item(32, 6).Value = 1615
$ws.Cells.Item(33, 6).Value = 1011
$ws.Cells.Item(34, 6).Value = 118
$ws.Cells.Item(35, 6).Value = 18
$ws.Cells.Item(36, 6).Value = 1107
$ws.Cells.Item(37, 6).Value = 2039
$ws.Cells.Item(38, 6).Value = 260
$ws.Cells.Item(41, 6).Value = 85
$ws.Cells.Item(42, 6).Value = 23
$ws.Cells.Item(43, 6).Value = 653
$ws.Cells.Item(44, 6).Value = 1312
$ws.Cells.Item(45, 6).Value = 89
$ws.Cells.Item(46, 6).Value = 79
$ws.Cells.Item(47, 6).Value = 428
$ws.Cells.Item(48, 6).Value = 64

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 10
$ws.Cells.Item(5, 6).Value = 66
$ws.Cells.Item(13, 6).Value = 12

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 883
$ws.Cells.Item(3, 6).Value = 1185
$ws.Cells.Item(4, 6).Value = 10
$ws.Cells.Item(5, 6).Value = 67
$ws.Cells.Item(6, 6).Value = 4346
$ws.Cells.Item(7, 6).Value = 2592
$ws.Cells.Item(8, 6).Value = 2505
$ws.Cells.Item(9, 6).Value = 1655
$ws.Cells.Item(12, 6).Value = 658
$ws.Cells.Item(13, 6).Value = 19
$ws.Cells.Item(15, 6).Value = 319
$ws.Cells.Item(17, 6).Value = 271
$ws.Cells.Item(18, 6).Value = 74
$ws.Cells.Item(19, 6).Value = 476
$ws.Cells.Item(22, 6).Value = 540
$ws.Cells.Item(24, 6).Value = 102
$ws.Cells.Item(25, 6).Value = 66
$ws.Cells.Item(29, 6).Value = 400
$ws.Cells.Item(30, 6).Value = 1615
$ws.Cells.Item(31, 6).Value = 1011
$ws.Cells.Item(32, 6).Value = 118
$ws.Cells.Item(34, 6).Value = 2039
$ws.Cells.Item(35, 6).Value = 260
$ws.Cells.Item(41, 6).Value = 85
$ws.Cells.Item(42, 6).Value = 23
$ws.Cells.Item(43, 6).Value = 653
$ws.Cells.Item(44, 6).Value = 1312
$ws.Cells.Item(46, 6).Value = 89
$ws.Cells.Item(47, 6).Value = 428
$ws.Cells.Item(48, 6).Value = 64
$ws.Cells.Item(49, 6).Value = 12
